$d = $word.ActiveDocument

# 1. In the last paragraph ("Jason : Fichier du menu ...") insert "(" (with a
#    leading space) right before "+ Ajouter et supprimer plat ...".
$d.Content.Find.Execute(
    "Jason : Fichier du menu + Ajouter et supprimer plat + les afficher",
    $true, $false, $false, $false, $false, $true, 1, $false,
    "Jason : Fichier du menu ( + Ajouter et supprimer plat + les afficher",
    2
)

# 2. Append two new paragraphs at the very end of the document:
#      - an empty paragraph
#      - a paragraph containing the text "test"
$endRange = $d.Content
$endRange.Collapse(0)          # wdCollapseEnd
$endRange.InsertParagraphAfter()
$endRange.Collapse(0)
$endRange.InsertParagraphAfter()

$finalRange = $d.Content
$finalRange.Collapse(0)
$finalRange.InsertAfter("test")

$d.Save()
